$d = $word.ActiveDocument

# 1. Global role rename: "Program Administrator" -> "Program Manager"
#    (also fixes the plural "Program Administrators" -> "Program Managers"
#    since the singular text is a prefix of the plural text)
$d.Content.Find.Execute("Program Administrator", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Program Manager", 2)

# 2. Scenario 1 intro paragraph: add a new closing sentence about the Lecturer
$d.Content.Find.Execute("(Art, Computer Science, Economics). ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "(Art, Computer Science, Economics). The Lecturer has not yet submitted any constraint information.", 2)

# 3. "This particular Program Manager is in charge of the Computer Science department " -> trailing space becomes a period
$d.Content.Find.Execute("This particular Program Manager is in charge of the Computer Science department ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "This particular Program Manager is in charge of the Computer Science department.", 2)

# 4. "Once logged in the Program Manager opens ... department" -> add comma after "in" and period at end
$d.Content.Find.Execute("Once logged in the Program Manager opens a database of all listed courses for the Computer Science department", $true, $false, $false, $false, $false,
                         $true, 1, $false, "Once logged in, the Program Manager opens a database of all listed courses for the Computer Science department.", 2)

# 5. "This database will be a list of courses (in this case Computer Science courses) offered at UCSC" -> add trailing period
$d.Content.Find.Execute("This database will be a list of courses (in this case Computer Science courses) offered at UCSC", $true, $false, $false, $false, $false,
                         $true, 1, $false, "This database will be a list of courses (in this case Computer Science courses) offered at UCSC.", 2)

# 6. "The database will include course context, teacher (empty value), time (empty value), location (empty value)" -> add trailing period
$d.Content.Find.Execute("The database will include course context, teacher (empty value), time (empty value), location (empty value)", $true, $false, $false, $false, $false,
                         $true, 1, $false, "The database will include course context, teacher (empty value), time (empty value), location (empty value).", 2)

# 7. "The Program Manager selects the courses that will be offered for the particular quarter" -> "the" becomes "this", add trailing period
$d.Content.Find.Execute("The Program Manager selects the courses that will be offered for the particular quarter", $true, $false, $false, $false, $false,
                         $true, 1, $false, "The Program Manager selects the courses that will be offered for this particular quarter.", 2)

# 8. Final bullet: append the new clause about the Program Administrator's final approval
$d.Content.Find.Execute("titled: “Fall 2010” ", $true, $false, $false, $false, $false,
                         $true, 1, $false, "titled: “Fall 2010” where they will be reviewed by the Program Administrator for final approval.", 2)
